$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B10").Value = "'0.21"
$ws.Range("C10").Value = "'0.42"
$ws.Range("D10").Value = "'0.3"
$ws.Range("E10").Value = "'0.45"
$ws.Range("F10").Value = "'0.39"
$ws.Range("G10").Value = "'0.53"
$ws.Range("H10").Value = "'0.47"
$ws.Range("I10").Value = "'0.53"
$ws.Range("J10").Value = "'0.56"
$ws.Range("K10").Value = "'0.62"
$ws.Range("L10").Value = "'0.59"
